# NIT-9011159481.xlsx : "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# - Swap the "Novedad de Retiro" / "Novedad de Ingreso" header columns
# - Update Valor Mora total (E11) and the worker/period counters (C13 / F13)
# - Replace worker #1 data (row 16) with a new person (ANIBAL ANDRES MENDOZA GARCIA)
# - Insert 6 extra detail rows (18-23) with new workers/periods, shifting the old
#   signature block down from rows 23/24 to rows 29/30
# - Last detail row (24) becomes CAMILO PARDO LONDOÃ?O / period 2506

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header summary block ------------------------------------------------
$ws.Range("E11").Value = 675157
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 7

# Swap "Novedad de Retiro" / "Novedad de Ingreso" in the table header row
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# ---- Row 16: replace with new worker --------------------------------------
$ws.Range("C16").Value = "1047445080"
$ws.Range("D16").Value = "ANIBAL ANDRES MENDOZA GARCIA"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 2378
$ws.Range("G16").Value = 1783800

# ---- Row 17: same worker as before, different period -----------------------
$ws.Range("E17").Value = "2302"

# ---- Insert 6 new detail rows below row 17 (rows 18-23), copying the
#      existing "middle" row formatting (row 17) so styles line up -----------
$ws.Rows("18:23").Insert()
$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 18
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1043995152"
$ws.Range("D18").Value = "ELIENETH RAMIREZ ESCOBAR"
$ws.Range("E18").Value = "2301"
$ws.Range("F18").Value = 54856
$ws.Range("G18").Value = 1551442

# Row 19
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1043995152"
$ws.Range("D19").Value = "ELIENETH RAMIREZ ESCOBAR"
$ws.Range("E19").Value = "2212"
$ws.Range("F19").Value = 47542
$ws.Range("G19").Value = 1551442

# Row 20
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1047419167"
$ws.Range("D20").Value = "LISETH TORRES VALENCIA"
$ws.Range("E20").Value = "1812"
$ws.Range("F20").Value = 1333
$ws.Range("G20").Value = 1000000

# Row 21
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1002189287"
$ws.Range("D21").Value = "RICHARD SAMITH MEDRANO MATOREL"
$ws.Range("E21").Value = "2507"
$ws.Range("F21").Value = 71352
$ws.Range("G21").Value = 1783800

# Row 22
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1002244395"
$ws.Range("D22").Value = "JIMMY JUNIOR SALAS CORTES"
$ws.Range("E22").Value = "2101"
$ws.Range("F22").Value = 42840
$ws.Range("G22").Value = 1071008

# Row 23
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1016105126"
$ws.Range("D23").Value = "CAMILO PARDO LONDOÃ?O"
$ws.Range("E23").Value = "2507"
$ws.Range("F23").Value = 200000
$ws.Range("G23").Value = 5000000

# ---- Row 24 (last detail row, keeps its distinct "last row" style) --------
$ws.Range("C24").Value = "1016105126"
$ws.Range("D24").Value = "CAMILO PARDO LONDOÃ?O"
$ws.Range("E24").Value = "2506"
$ws.Range("F24").Value = 200000
$ws.Range("G24").Value = 5000000

Write-Host "done"
